$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.98943825577132
$ws.Range("C2").Value = 8.454860364875108
$ws.Range("D2").Value = 7.826453265643561
$ws.Range("E2").Value = 13.18243162950818
$ws.Range("F2").Value = 39.21268724597822
$ws.Range("J2").Value = 10.46998192983836
$ws.Range("K2").Value = 11.17209005334947
$ws.Range("L2").Value = 10.75886745923918
$ws.Range("M2").Value = 15.85931230122848
$ws.Range("N2").Value = 22.14080014276757
$ws.Range("O2").Value = 30.35078949967882
$ws.Range("B3").Value = 14.83767363624774
$ws.Range("C3").Value = 8.433782680159677
$ws.Range("D3").Value = 7.818760874871302
$ws.Range("E3").Value = 13.20497416956307
$ws.Range("F3").Value = 39.29629240620887
$ws.Range("J3").Value = 10.48674579404068
$ws.Range("K3").Value = 11.05898004020352
$ws.Range("L3").Value = 10.76701763653766
$ws.Range("M3").Value = 15.84351969159543
$ws.Range("N3").Value = 22.2010166238085
$ws.Range("O3").Value = 30.4301639208446
$ws.Range("B4").Value = 14.74646009688573
$ws.Range("C4").Value = 8.4207100879427
$ws.Range("D4").Value = 7.815012378423155
$ws.Range("E4").Value = 13.22007704451477
$ws.Range("F4").Value = 39.3544338406386
$ws.Range("J4").Value = 10.49760597704258
$ws.Range("K4").Value = 10.99079257564639
$ws.Range("L4").Value = 10.77309331754179
$ws.Range("M4").Value = 15.8358324600591
$ws.Range("N4").Value = 22.23974107025943
$ws.Range("O4").Value = 30.48369821669938
$ws.Range("B5").Value = 14.70982404545595
$ws.Range("C5").Value = 8.415351013580509
$ws.Range("D5").Value = 7.813731502843487
$ws.Range("E5").Value = 13.22654939544054
$ws.Range("F5").Value = 39.37983782280898
$ws.Range("J5").Value = 10.50217458702363
$ws.Range("K5").Value = 10.96335117137787
$ws.Range("L5").Value = 10.77583929229319
$ws.Range("M5").Value = 15.83320834543768
$ws.Range("N5").Value = 22.25596323037436
$ws.Range("O5").Value = 30.50672002068924
$ws.Range("B6").Value = 14.70377403356782
$ws.Range("C6").Value = 8.414459246676701
$ws.Range("D6").Value = 7.813533756555787
$ws.Range("E6").Value = 13.22764333519221
$ws.Range("F6").Value = 39.38415943619483
$ws.Range("J6").Value = 10.50294184953439
$ws.Range("K6").Value = 10.95881622507057
$ws.Range("L6").Value = 10.77631159062848
$ws.Range("M6").Value = 15.83280341377083
$ws.Range("N6").Value = 22.25868361835613
$ws.Range("O6").Value = 30.51061560999858
$ws.Range("B7").Value = 14.74596379689704
$ws.Range("C7").Value = 8.420637941376651
$ws.Range("D7").Value = 7.814994103327525
$ws.Range("E7").Value = 13.2201630453829
$ws.Range("F7").Value = 39.35476952174234
$ws.Range("J7").Value = 10.49766701139809
$ws.Range("K7").Value = 10.99042105556567
$ws.Range("L7").Value = 10.77312925629354
$ws.Range("M7").Value = 15.83579500741014
$ws.Range("N7").Value = 22.2399580581567
$ws.Range("O7").Value = 30.48400381370375
$ws.Range("B8").Value = 14.93672432002265
$ws.Range("C8").Value = 8.447620531385494
$ws.Range("D8").Value = 7.823599616680181
$ws.Range("E8").Value = 13.18994277144606
$ws.Range("F8").Value = 39.24010089852858
$ws.Range("J8").Value = 10.47564464598139
$ws.Range("K8").Value = 11.13284469143176
$ws.Range("L8").Value = 10.76145566860045
$ws.Range("M8").Value = 15.85345191713099
$ws.Range("N8").Value = 22.16120006803234
$ws.Range("O8").Value = 30.37716174986309
$ws.Range("B9").Value = 15.32470085877999
$ws.Range("C9").Value = 8.499452317396232
$ws.Range("D9").Value = 7.848135364505056
$ws.Range("E9").Value = 13.14066723275279
$ws.Range("F9").Value = 39.06928683177451
$ws.Range("J9").Value = 10.43694055117991
$ws.Range("K9").Value = 11.42089747315581
$ws.Range("L9").Value = 10.74703558540002
$ws.Range("M9").Value = 15.90387347449309
$ws.Range("N9").Value = 22.02059308063252
$ws.Range("O9").Value = 30.20573099641377
$ws.Range("B10").Value = 15.61584971034157
$ws.Range("C10").Value = 8.536817980687825
$ws.Range("D10").Value = 7.870728217622279
$ws.Range("E10").Value = 13.11051928720114
$ws.Range("F10").Value = 38.97678120092618
$ws.Range("J10").Value = 10.41121161155599
$ws.Range("K10").Value = 11.63616428011431
$ws.Range("L10").Value = 10.74156529659041
$ws.Range("M10").Value = 15.95033816151686
$ws.Range("N10").Value = 21.92564140667564
$ws.Range("O10").Value = 30.10301880039655
$ws.Range("B11").Value = 15.7491178999304
$ws.Range("C11").Value = 8.553649251027885
$ws.Range("D11").Value = 7.881973576386063
$ws.Range("E11").Value = 13.09811181105012
$ws.Range("F11").Value = 38.94186606478449
$ws.Range("J11").Value = 10.4000892928218
$ws.Range("K11").Value = 11.73451927196706
$ws.Range("L11").Value = 10.74018015489201
$ws.Range("M11").Value = 15.97347222784692
$ws.Range("N11").Value = 21.88424128782863
$ws.Range("O11").Value = 30.06134145885405
$ws.Range("B12").Value = 15.79965906265242
$ws.Range("C12").Value = 8.559997775706563
$ws.Range("D12").Value = 7.886368809445937
$ws.Range("E12").Value = 13.09360077864878
$ws.Range("F12").Value = 38.92967499374423
$ws.Range("J12").Value = 10.39596083502762
$ws.Range("K12").Value = 11.77179507604473
$ws.Range("L12").Value = 10.73981344650792
$ws.Range("M12").Value = 15.98251507315575
$ws.Range("N12").Value = 21.8688207883532
$ws.Range("O12").Value = 30.04628526611083
$ws.Range("B13").Value = 15.78877144490823
$ws.Range("C13").Value = 8.558631642775449
$ws.Range("D13").Value = 7.885416168297475
$ws.Range("E13").Value = 13.0945639836473
$ws.Range("F13").Value = 38.93225472780937
$ws.Range("J13").Value = 10.39684627207723
$ws.Range("K13").Value = 11.76376617638724
$ws.Range("L13").Value = 10.73988541865332
$ws.Range("M13").Value = 15.98055504857865
$ws.Range("N13").Value = 21.87213046424669
$ws.Range("O13").Value = 30.04949558752908
$ws.Range("B14").Value = 15.75327468526613
$ws.Range("C14").Value = 8.554172057939205
$ws.Range("D14").Value = 7.88233244555116
$ws.Range("E14").Value = 13.09773693280484
$ws.Range("F14").Value = 38.94084244526817
$ws.Range("J14").Value = 10.39974797421605
$ws.Range("K14").Value = 11.73758553045734
$ws.Range("L14").Value = 10.74014682767893
$ws.Range("M14").Value = 15.97421055462437
$ws.Range("N14").Value = 21.88296749329284
$ws.Range("O14").Value = 30.06008822261166
$ws.Range("B15").Value = 15.73154044598958
$ws.Range("C15").Value = 8.551437128511242
$ws.Range("D15").Value = 7.880461329338585
$ws.Range("E15").Value = 13.09970484612815
$ws.Range("F15").Value = 38.94623687298803
$ws.Range("J15").Value = 10.40153619120933
$ws.Range("K15").Value = 11.72155224599953
$ws.Range("L15").Value = 10.74032747479933
$ws.Range("M15").Value = 15.97036101242619
$ws.Range("N15").Value = 21.88963889961046
$ws.Range("O15").Value = 30.06667108829571
$ws.Range("B16").Value = 15.6071532005185
$ws.Range("C16").Value = 8.535714566944696
$ws.Range("D16").Value = 7.870012586366229
$ws.Range("E16").Value = 13.11135639932392
$ws.Range("F16").Value = 38.97920719867742
$ws.Range("J16").Value = 10.41195016029251
$ws.Range("K16").Value = 11.62974254308157
$ws.Range("L16").Value = 10.74167794938631
$ws.Range("M16").Value = 15.9488660852755
$ws.Range("N16").Value = 21.92838300230257
$ws.Range("O16").Value = 30.10584410324226
$ws.Range("B17").Value = 15.531025767671
$ws.Range("C17").Value = 8.526025801755628
$ws.Range("D17").Value = 7.863848882916423
$ws.Range("E17").Value = 13.11883862131261
$ws.Range("F17").Value = 39.00126896635681
$ws.Range("J17").Value = 10.41848757964256
$ws.Range("K17").Value = 11.57350813464904
$ws.Range("L17").Value = 10.7427884214318
$ws.Range("M17").Value = 15.93618786601277
$ws.Range("N17").Value = 21.95260990114662
$ws.Range("O17").Value = 30.13116840822898
$ws.Range("B18").Value = 15.48731930791014
$ws.Range("C18").Value = 8.520437473154905
$ws.Range("D18").Value = 7.860394918275672
$ws.Range("E18").Value = 13.12326524847917
$ws.Range("F18").Value = 39.01463283611648
$ws.Range("J18").Value = 10.42230252295592
$ws.Range("K18").Value = 11.54120579832069
$ws.Range("L18").Value = 10.7435310379164
$ws.Range("M18").Value = 15.92908395670647
$ws.Range("N18").Value = 21.96671348731489
$ws.Range("O18").Value = 30.14620930418271
$ws.Range("B19").Value = 15.47253611254309
$ws.Range("C19").Value = 8.518542715817725
$ws.Range("D19").Value = 7.859241203193368
$ws.Range("E19").Value = 13.12478517859273
$ws.Range("F19").Value = 39.01927345713802
$ws.Range("J19").Value = 10.42360361874854
$ws.Range("K19").Value = 11.53027697551959
$ws.Range("L19").Value = 10.74380034546471
$ws.Range("M19").Value = 15.92671116951384
$ws.Range("N19").Value = 21.97151776514185
$ws.Range("O19").Value = 30.15138346686859
$ws.Range("B20").Value = 15.53912169315856
$ws.Range("C20").Value = 8.527058811017977
$ws.Range("D20").Value = 7.864495593578555
$ws.Range("E20").Value = 13.11802939452996
$ws.Range("F20").Value = 38.99885064154706
$ws.Range("J20").Value = 10.41778599194037
$ws.Range("K20").Value = 11.5794902392584
$ws.Range("L20").Value = 10.74265946226693
$ws.Range("M20").Value = 15.93751803246243
$ws.Range("N20").Value = 21.95001343257706
$ws.Range("O20").Value = 30.12842342647603
$ws.Range("B21").Value = 15.76369923798956
$ws.Range("C21").Value = 8.555482636939516
$ws.Range("D21").Value = 7.883234514220745
$ws.Range("E21").Value = 13.09679987875776
$ws.Range("F21").Value = 38.93829205702276
$ws.Range("J21").Value = 10.39889341540593
$ws.Range("K21").Value = 11.74527482557365
$ws.Range("L21").Value = 10.74006576935139
$ws.Range("M21").Value = 15.97606645850604
$ws.Range("N21").Value = 21.87977743390694
$ws.Range("O21").Value = 30.05695719989134
$ws.Range("B22").Value = 15.91089353748763
$ws.Range("C22").Value = 8.573912333319793
$ws.Range("D22").Value = 7.896278062267606
$ws.Range("E22").Value = 13.084017277325
$ws.Range("F22").Value = 38.90472009127937
$ws.Range("J22").Value = 10.38703153230924
$ws.Range("K22").Value = 11.85379078331015
$ws.Range("L22").Value = 10.73929008997693
$ws.Range("M22").Value = 16.00290423633161
$ws.Range("N22").Value = 21.83537061650434
$ws.Range("O22").Value = 30.0144823433642
$ws.Range("B23").Value = 15.83230881152288
$ws.Range("C23").Value = 8.564089886189681
$ws.Range("D23").Value = 7.889244376313422
$ws.Range("E23").Value = 13.09073984131761
$ws.Range("F23").Value = 38.92208853217473
$ws.Range("J23").Value = 10.39331813472499
$ws.Range("K23").Value = 11.79586863825745
$ws.Range("L23").Value = 10.73962024239197
$ws.Range("M23").Value = 15.98843157849131
$ws.Range("N23").Value = 21.85893481425629
$ws.Range("O23").Value = 30.03676460190017
$ws.Range("B24").Value = 15.53546133439729
$ws.Range("C24").Value = 8.526591843916858
$ws.Range("D24").Value = 7.864202936336789
$ws.Range("E24").Value = 13.11839485656883
$ws.Range("F24").Value = 38.99994184722024
$ws.Range("J24").Value = 10.41810300374071
$ws.Range("K24").Value = 11.57678564151819
$ws.Range("L24").Value = 10.74271744010079
$ws.Range("M24").Value = 15.9369160877023
$ws.Range("N24").Value = 21.95118675012449
$ws.Range("O24").Value = 30.12966293254657
$ws.Range("B25").Value = 15.2184985660242
$ws.Range("C25").Value = 8.485551925347144
$ws.Range("D25").Value = 7.84068793842849
$ws.Range("E25").Value = 13.1529318681037
$ws.Range("F25").Value = 39.10970493415697
$ws.Range("J25").Value = 10.44693385594804
$ws.Range("K25").Value = 11.34220918585575
$ws.Range("L25").Value = 10.75003368283882
$ws.Range("M25").Value = 15.88856272663745
$ws.Range("N25").Value = 22.05715827457752
$ws.Range("O25").Value = 30.24802790656679
